# Applies a row-wise data shuffle (weekly refresh of source rows) to the
# "Hortaliza, Mapocho Venta Directa de Santiago - Ají" sheet.
#
# The columns Fecha(D), Variedad(H), Calidad(I), Volumen(J), Precio minimo(K),
# Precio maximo(L), Precio promedio ponderado(M), Unidad de comercializacion(N),
# Precio $/Kg(P) and Kg o Unidades(Q) for data rows 2-18 get rearranged: the
# full set of those column values that used to live on one row now lives on
# another row (a permutation of the 17 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values travel together as a group per source row.
$cols = @("D", "H", "I", "J", "K", "L", "M", "N", "P", "Q")

# Destination row -> source row (i.e. destRow ends up holding what srcRow held).
$mapping = @{
    2  = 10
    3  = 9
    4  = 11
    5  = 13
    6  = 3
    7  = 5
    8  = 17
    9  = 18
    10 = 6
    11 = 7
    12 = 14
    13 = 2
    14 = 4
    15 = 12
    16 = 15
    17 = 8
    18 = 16
}

# Snapshot the original values for every row/column involved before writing
# anything back, since several rows read from each other. The Fecha column
# (D) is date-formatted, so use Value2() to capture the underlying numeric
# date serial rather than a formatted date string.
$snapshot = @{}
foreach ($row in 2..18) {
    $rowValues = @{}
    foreach ($col in $cols) {
        if ($col -eq "D") {
            $rowValues[$col] = $ws.Range("$col$row").Value2()
        } else {
            $rowValues[$col] = $ws.Range("$col$row").Value()
        }
    }
    $snapshot[$row] = $rowValues
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcValues = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcValues[$col]
    }
}
